# Fix header label for "Requested quantity" -> new metric names,
# and add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast,
# yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

# 1) Weekly Quantity sheet: rename the quantity header.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Monthly Trend sheet: rename the quantity header.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add the new "PO Forecast" sheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the sheet-level layout settings used by the other sheets.
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Reuse the bold/centered header style and the date number-format style
# already present in the workbook (copy formats only, values are set below).
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A50").PasteSpecial(-4122)

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper).
$rows = @(
    @(44990.99999999999,0,-8.282791521422766,8.299787722190308),
    @(45004.99999999999,1,-6.977512095854269,8.364000941616704),
    @(45011.99999999999,1,-6.904683469599772,8.725668313846013),
    @(45025.99999999999,1,-6.957675821025259,8.377822209855854),
    @(45046.99999999999,2,-6.448975829041137,9.701206481364144),
    @(45060.99999999999,2,-6.448522130065042,9.669411121592622),
    @(45088.99999999999,3,-4.853013002698611,10.58590234418656),
    @(45095.99999999999,3,-5.498863021527043,10.79460309468888),
    @(45102.99999999999,3,-4.903029261639294,10.47630901526262),
    @(45109.99999999999,3,-5.26581417852655,11.20468916480163),
    @(45116.99999999999,3,-4.863755895605586,10.90091846883878),
    @(45130.99999999999,4,-5.06675690444781,11.54724841093604),
    @(45158.99999999999,4,-3.66563793049982,12.01463252220325),
    @(45172.99999999999,5,-3.38517749914846,12.87773661260587),
    @(45179.99999999999,5,-3.043467459374462,12.51000704864666),
    @(45200.99999999999,5,-2.955556907464479,13.32173240464704),
    @(45207.99999999999,5,-3.184019046900614,13.36604722725671),
    @(45228.99999999999,6,-2.208804616062694,14.09340467655948),
    @(45242.99999999999,6,-2.152516171627181,14.02947965195865),
    @(45256.99999999999,6,-1.297754917323064,14.01588491096251),
    @(45270.99999999999,7,-1.371535548947921,15.06245252973236),
    @(45277.99999999999,7,-1.322034490211108,14.87205355017154),
    @(45305.99999999999,8,-0.07208052474946607,15.50233518812161),
    @(45312.99999999999,8,-0.5865736239344982,15.23300036281984),
    @(45319.99999999999,8,0.3193233940193386,15.8346200506836),
    @(45333.99999999999,8,-0.05461583998993431,16.26310095015166),
    @(45340.99999999999,8,-0.1460721935359754,15.70341610727905),
    @(45354.99999999999,9,0.8115541225475902,16.62359801022781),
    @(45361.99999999999,9,1.351651769293508,16.97076910509673),
    @(45368.99999999999,9,0.6684691882658738,17.19916915869037),
    @(45375.99999999999,9,1.205090968111695,17.41711911967925),
    @(45382.99999999999,9,1.194561638917558,17.66284630802485),
    @(45389.99999999999,9,1.883545588691104,17.81120741494591),
    @(45396.99999999999,10,1.143024022404023,17.10516261202403),
    @(45438.99999999999,11,2.416530569840205,18.76788680642683),
    @(45522.99999999999,13,3.77847481092258,19.72508200295777),
    @(45536.99999999999,13,4.703698355738106,20.80650485969697),
    @(45543.99999999999,13,4.693160840461046,20.68699455524201),
    @(45550.99999999999,13,5.00024663521881,20.99170831388316),
    @(45571.99999999999,14,5.817442180644965,21.96285306386334),
    @(45578.99999999999,14,6.240709463361801,21.64150681546626),
    @(45585.99999999999,14,5.454440284680373,21.63804875989485),
    @(45592.99999999999,14,6.73081472143606,22.10654146956426),
    @(45599.99999999999,14,6.693873404053464,23.12439871126898),
    @(45606.99999999999,14,6.417422185588268,22.31254467354906),
    @(45613.99999999999,15,6.984892873842216,23.12979573317048),
    @(45620.99999999999,15,6.913689885142341,22.7021822285467),
    @(45627.99999999999,15,6.954632314876347,22.93566595002),
    @(45634.99999999999,15,7.01016568602047,23.99506238076033)
)

$n = $rows.Count
$data = New-Object 'object[,]' $n,4
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}
$wsForecast.Range("A2:D50").Value = $data

# Keep the originally active sheet selected (unchanged by this edit).
$null = $wsWeekly.Activate()
$null = $wsWeekly.Range("A1").Select()
